# Update auf 2021 Daten - finales Update vor Repo Migration
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column E header: "Abschlussjahr Doppelter Jahrgang" -> "Jahr mit weniger Abschlüssen"
$ws.Range("E1").Value = "Jahr mit weniger Abschlüssen"

# Update column E values (Jahr mit weniger Abschlüssen) with 2021 data
$ws.Range("E2").Value = 2020
$ws.Range("E3").Value = "2025, 2026"
$ws.Range("E8").Value = "-"
$ws.Range("E10").Value = 2023
$ws.Range("E11").Value = 2027
$ws.Range("E16").Value = 2027

# Selection moved to E17 as last active cell of the editing session
$ws.Range("E17").Select()
